$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warheads")

# Update sheet view: remove topLeftCell scrolling, set new selection to N12
$ws.Activate()
$ws.Range("N12").Select()

# Update N column values per diff
$ws.Range("N2").Value = 0.1
$ws.Range("N3").Value = 0.1
$ws.Range("N4").Value = 0.125
$ws.Range("N6").Value = 0.1
$ws.Range("N7").Value = 0.1
$ws.Range("N8").Value = 0.2
$ws.Range("N10").Value = 0.1
$ws.Range("N11").Value = 0.1
$ws.Range("N12").Value = 0.2
$ws.Range("N20").Value = 0.1
$ws.Range("N21").Value = 0.1
$ws.Range("N22").Value = 0.3
$ws.Range("N24").Value = 0.1
$ws.Range("N26").Value = 0.4
$ws.Range("N28").Value = 0.1
$ws.Range("N29").Value = 0.1
$ws.Range("N30").Value = 0.4
$ws.Range("N38").Value = 0.05
$ws.Range("N39").Value = 0.05
$ws.Range("N40").Value = 0.6
$ws.Range("N42").Value = 0.025
$ws.Range("N43").Value = 0.025
$ws.Range("N46").Value = 0.025
$ws.Range("N47").Value = 0.025
